# Updates cryptos list prices/volumes per the latest scrape.
# Values in column D that look numeric are written with a leading
# apostrophe so Excel stores them as literal text (matching the
# workbook convention of keeping price strings, e.g. "237.30", as text)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.211.71'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '1.829.95'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("D5").Value = '''237.30'
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").Value = '''0.6091'
$ws.Range("E6").Value = '  -3.94%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = '''0.07077'
$ws.Range("E8").Value = '  -5.34%  '
$ws.Range("D9").Value = '''0.2816'
$ws.Range("E9").Value = '  -3.11%  '
$ws.Range("D10").Value = '''23.78'
$ws.Range("E10").Value = '  -5.53%  '
$ws.Range("D11").Value = '''0.07643'
$ws.Range("E11").Value = '  -1.34%  '
$ws.Range("D12").Value = '1.830.94'
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("D13").Value = '''4.810'
$ws.Range("E13").Value = '  -3.67%  '
$ws.Range("D14").Value = '''0.6324'
$ws.Range("E14").Value = '  -7.04%  '
$ws.Range("D15").Value = '''0.000009971'
$ws.Range("E15").Value = '  -2.64%  '
$ws.Range("D16").Value = '2.073.21'
$ws.Range("E16").Value = '  -1.30%  '
$ws.Range("D17").Value = '''79.40'
$ws.Range("E17").Value = '  -3.24%  '
$ws.Range("D18").Value = '''5.952'
$ws.Range("E18").Value = '  -5.12%  '
$ws.Range("D19").Value = '29.224.03'
$ws.Range("E19").Value = '  -0.51%  '
$ws.Range("D20").Value = '''228.60'
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '''1.002'
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").Value = '''11.80'
$ws.Range("E22").Value = '  -4.43%  '
$ws.Range("D23").Value = '''7.033'
$ws.Range("E23").Value = '  -5.40%  '
$ws.Range("D24").Value = '''1.002'
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").Value = '''155.50'
$ws.Range("E25").Value = '  -1.70%  '
$ws.Range("D26").Value = '''8.098'
$ws.Range("E26").Value = '  -4.86%  '
$ws.Range("D27").Value = '''0.1301'
$ws.Range("E27").Value = '  -4.37%  '
$ws.Range("D28").Value = '''16.72'
$ws.Range("E28").Value = '  -4.56%  '
$ws.Range("D29").Value = '''0.06743'
$ws.Range("E29").Value = '  +3.16%  '
$ws.Range("D30").Value = '''1.480'
$ws.Range("E30").Value = '  +3.36%  '
$ws.Range("E31").Value = '  -1.87%  '
$ws.Range("D32").Value = '''3.846'
$ws.Range("E32").Value = '  -5.72%  '
$ws.Range("D33").Value = '''3.831'
$ws.Range("E33").Value = '  -5.50%  '
$ws.Range("D34").Value = '''1.131'
$ws.Range("E34").Value = '  -1.04%  '
$ws.Range("D35").Value = '''1.737'
$ws.Range("E35").Value = '  -5.73%  '
$ws.Range("D36").Value = '''0.6547'
$ws.Range("E36").Value = '  -6.51%  '
$ws.Range("D37").Value = '''2.559'
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("D38").Value = '1.236.21'
$ws.Range("E38").Value = '  -1.19%  '
$ws.Range("D39").Value = '''2.762'
$ws.Range("E39").Value = '  -1.99%  '
$ws.Range("D40").Value = '''0.01766'
$ws.Range("D41").Value = '''6.585'
$ws.Range("E41").Value = '  -2.79%  '
$ws.Range("D42").Value = '''0.9227'
$ws.Range("E42").Value = '  -1.39%  '
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").Value = '1.986.30'
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("D45").Value = '''100.93'
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("D46").Value = '''63.52'
$ws.Range("E46").Value = '  -3.00%  '
$ws.Range("D47").Value = '''0.00000000116'
$ws.Range("D48").Value = '''1.627'
$ws.Range("E48").Value = '  -5.50%  '
$ws.Range("D49").Value = '''8.574'
$ws.Range("E49").Value = '  -5.36%  '
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").Value = '''6.522'
$ws.Range("E50").Value = '  -7.80%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '''0.1083'
$ws.Range("E51").Value = '  -5.77%  '
